$wb = $excel.ActiveWorkbook

# Row 40 values per sheet: B, C, D, E (text/hex-like strings), F, G, H, I (numbers)
# A40 is the same timestamp for every sheet.

$timeValue = 45826.43626157408

function Add-Row40 {
    param($ws, $bVal, $cVal, $dVal, $eVal, $fVal, $gStr, $hVal, $iVal)

    $ws.Range("A40").Value = $timeValue
    $ws.Range("A40").NumberFormat = $ws.Range("A39").NumberFormat

    $ws.Range("B40").Value = $bVal
    $ws.Range("C40").Value = $cVal
    $ws.Range("D40").Value = $dVal
    $ws.Range("E40").Value = $eVal

    $ws.Range("F40").Value = $fVal
    # Go through numeric (not string) arithmetic so the cell stays a plain
    # number (no auto-applied numeric format / style) while keeping full
    # double precision of the literal.
    $gNum = 0 + $gStr
    $ws.Range("G40").Value = $gNum
    $ws.Range("H40").Value = $hVal
    $ws.Range("I40").Value = $iVal
}

$ws1 = $wb.Worksheets.Item("DE_LFT_#1")
Add-Row40 $ws1 "0x01,0x7c" "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0," "0x01,0x6C" "0x14" 380 "7.598631275147109E+23" 364 14

$ws2 = $wb.Worksheets.Item("DE_LFT_#2")
Add-Row40 $ws2 "0x01,0x7c" "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78," "0x01,0x6C" "0xe" 380 "5.68432987514711E+23" 364 14

$ws3 = $wb.Worksheets.Item("DE_PLT_#1")
Add-Row40 $ws3 "0x00,0x82" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x00,0x80" "0x7" 130 "5.68631262647114E+23" 128 7

$ws4 = $wb.Worksheets.Item("DE_PLT_#2")
Add-Row40 $ws4 "0x00,0x82" "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c," "0x00,0x80" "0x3" 130 "9.85046333984776E+23" 128 3
